$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dotenvs")
$tbl = $ws.ListObjects.Item("Table1")

# Insert two new rows right after the current last table row (row 9), copying
# formatting from the row above so the new cells pick up the same styles
# (wrap text, centered boolean column, etc.) used throughout the table.
$ws.Rows.Item(10).Insert(-4121, 0) | Out-Null
$ws.Range("A10").Value = "Database"
$ws.Range("B10").Value = "db_type"
$ws.Range("C10").Value = "The type of database to use"
$ws.Range("D10").Value = "mongodb"

$ws.Rows.Item(11).Insert(-4121, 0) | Out-Null
$ws.Range("A11").Value = "Logging"
$ws.Range("B11").Value = "log_db"
$ws.Range("C11").Value = "Log database"
$ws.Range("D11").Value = $false

# Grow the table (ListObject) so the new rows become part of it.
$tbl.Resize($ws.Range("A2:D11"))

$ws.Range("G14").Select()
